# Updates the "Price" (D) and "Volume(1h)" (E) columns in the crypto
# tracker sheet with freshly scraped values, matching the automated
# "Updated symbol list" GitHub Actions commit.
# Values are written as text (leading apostrophe forces text, matching
# the original inline-string cells), then the style is reset to "Normal"
# so no stray number formatting is introduced on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB
$ws.Range("D2").Value = "'307.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.12%"
$ws.Range("E2").Style = "Normal"

# Row 3: OKB
$ws.Range("D3").Value = "'36.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.91%"
$ws.Range("E3").Style = "Normal"

# Row 4: HuobiToken
$ws.Range("D4").Value = "'5.057"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.06%"
$ws.Range("E4").Style = "Normal"

# Row 5: Cronos
$ws.Range("D5").Value = "'0.08111"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.25%"
$ws.Range("E5").Style = "Normal"

# Row 6: FTXToken
$ws.Range("D6").Value = "'1.984"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'2.09%"
$ws.Range("E6").Style = "Normal"

# Row 7: GateToken
$ws.Range("D7").Value = "'4.165"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.50%"
$ws.Range("E7").Style = "Normal"

# Row 8: KuCoinToken
$ws.Range("D8").Value = "'7.880"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.02%"
$ws.Range("E8").Style = "Normal"

# Row 9: MXToken
$ws.Range("D9").Value = "'0.9294"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.22%"
$ws.Range("E9").Style = "Normal"

# Row 10: LiechtensteinCryptoassetsExchange
$ws.Range("D10").Value = "'0.1458"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'16.28%"
$ws.Range("E10").Style = "Normal"

# Row 11: WazirX
$ws.Range("D11").Value = "'0.1922"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.99%"
$ws.Range("E11").Style = "Normal"

# Row 12: MandalaExchangeToken
$ws.Range("D12").Value = "'0.09139"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.78%"
$ws.Range("E12").Style = "Normal"

# Row 13: BitrueCoin
$ws.Range("D13").Value = "'0.03442"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-2.13%"
$ws.Range("E13").Style = "Normal"

# Row 14: BitMartToken
$ws.Range("D14").Value = "'0.09886"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.42%"
$ws.Range("E14").Style = "Normal"

# Row 15: BitForexToken
$ws.Range("D15").Value = "'0.001425"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.19%"
$ws.Range("E15").Style = "Normal"

# Row 16: TigerCash
$ws.Range("D16").Value = "'0.006777"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'4.19%"
$ws.Range("E16").Style = "Normal"

# Row 17: LEO
$ws.Range("D17").Value = "'3.832"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'6.07%"
$ws.Range("E17").Style = "Normal"

# Row 18: BTSEToken
$ws.Range("E18").Value = "'10.42%"
$ws.Range("E18").Style = "Normal"

# Row 19: BitpandaEcosystemToken
$ws.Range("D19").Value = "'0.3453"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.44%"
$ws.Range("E19").Style = "Normal"

# Row 20: ProBitToken
$ws.Range("D20").Value = "'0.1314"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.13%"
$ws.Range("E20").Style = "Normal"

# Row 21: MCDex
$ws.Range("D21").Value = "'4.817"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-6.64%"
$ws.Range("E21").Style = "Normal"

# Row 23: CoinExToken
$ws.Range("D23").Value = "'0.04400"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.16%"
$ws.Range("E23").Style = "Normal"

# Row 24: BitKan
$ws.Range("D24").Value = "'0.001237"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.28%"
$ws.Range("E24").Style = "Normal"

# Row 25: HotbitToken
$ws.Range("D25").Value = "'0.004178"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-11.50%"
$ws.Range("E25").Style = "Normal"

# Row 27: NitroEx
$ws.Range("D27").Value = "'0.0001306"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'0.53%"
$ws.Range("E27").Style = "Normal"

# Row 39: One
$ws.Range("D39").Value = "'0.02036"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'3.78%"
$ws.Range("E39").Style = "Normal"

# Row 40: IDEX
$ws.Range("D40").Value = "'0.05134"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.30%"
$ws.Range("E40").Style = "Normal"

# Row 41: KickToken
$ws.Range("D41").Value = "'0.007482"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.91%"
$ws.Range("E41").Style = "Normal"

# Row 42: Dexo
$ws.Range("D42").Value = "'0.01015"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.07%"
$ws.Range("E42").Style = "Normal"

# Row 43: BKEXToken
$ws.Range("E43").Value = "'-0.12%"
$ws.Range("E43").Style = "Normal"

# Row 44: CEJI
$ws.Range("D44").Value = "'0.002129"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'1.50%"
$ws.Range("E44").Style = "Normal"

# Row 45: LocalTraders
$ws.Range("D45").Value = "'0.009888"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-7.69%"
$ws.Range("E45").Style = "Normal"

# Row 46: CoinLion
$ws.Range("D46").Value = "'0.00006311"
$ws.Range("D46").Style = "Normal"

# Row 47: Kangarootoken
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.19%"
$ws.Range("E47").Style = "Normal"

# Row 48: BOLO
$ws.Range("D48").Value = "'64.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-0.56%"
$ws.Range("E48").Style = "Normal"

# Row 49: CoinbaseStockToken
$ws.Range("D49").Value = "'0.001605"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-3.37%"
$ws.Range("E49").Style = "Normal"

# Row 50: CryptobidCoin
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.19%"
$ws.Range("E50").Style = "Normal"

# Row 51: SpecialPowerGold
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.19%"
$ws.Range("E51").Style = "Normal"
